# DDS_Project_2_Presentation.pptx - slide 21 ("Content Placeholder 3")
# The paragraph "We built a model to predict salary had an RMSE of 1080 and
# P-Value of < 2.2e-16" gets a trailing period added after "< 2.2e-16", which
# PowerPoint records by splitting the run right before the "<" and adding a
# new run for "< 2.2e-16." with the appended period.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(21)
$shape = $s.Shapes.Item(6)
$textRange = $shape.TextFrame.TextRange

$fullText = $textRange.Text
$target = "< 2.2e-16"
$startIdx = $fullText.IndexOf($target)

if ($startIdx -ge 0) {
    $oldRun = $textRange.Characters($startIdx + 1, $target.Length)
    $oldRun.Text = "< 2.2e-16."
}
